$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 170-171; everything previously at 170.. shifts down to 172..
$ws.Rows("170:171").Insert()

# New row 170: Black Amber / Primera, fecha 2022-03-07 (serial 44627), O'Higgins
$ws.Range("A170").Value = 8
$ws.Range("B170").Value = "Terminal La Palmera de La Serena"
$ws.Range("C170").Value = "Coquimbo"
$ws.Range("D170").Value = 44627
$ws.Range("E170").Value = 4
$ws.Range("F170").Value = "Fruta"
$ws.Range("G170").Value = 100103
$ws.Range("H170").Value = "Frutos de hueso (carozo)"
$ws.Range("I170").Value = 100103002
$ws.Range("J170").Value = "Ciruela"
$ws.Range("K170").Value = "Black Amber"
$ws.Range("L170").Value = "Primera"
$ws.Range("M170").Value = 20
$ws.Range("N170").Value = 220000
$ws.Range("O170").Value = 230000
$ws.Range("P170").Value = 225000
$ws.Range("Q170").Value = "`$/bins (450 kilos)"
$ws.Range("R170").Value = "Región de O'Higgins"
$ws.Range("S170").Value = 500
$ws.Range("T170").Value = 450

# New row 171: Larry Ann / Primera, fecha 2022-03-07 (serial 44627), O'Higgins
$ws.Range("A171").Value = 8
$ws.Range("B171").Value = "Terminal La Palmera de La Serena"
$ws.Range("C171").Value = "Coquimbo"
$ws.Range("D171").Value = 44627
$ws.Range("E171").Value = 4
$ws.Range("F171").Value = "Fruta"
$ws.Range("G171").Value = 100103
$ws.Range("H171").Value = "Frutos de hueso (carozo)"
$ws.Range("I171").Value = 100103002
$ws.Range("J171").Value = "Ciruela"
$ws.Range("K171").Value = "Larry Ann"
$ws.Range("L171").Value = "Primera"
$ws.Range("M171").Value = 20
$ws.Range("N171").Value = 220000
$ws.Range("O171").Value = 230000
$ws.Range("P171").Value = 225000
$ws.Range("Q171").Value = "`$/bins (450 kilos)"
$ws.Range("R171").Value = "Región de O'Higgins"
$ws.Range("S171").Value = 500
$ws.Range("T171").Value = 450

Write-Host "Rows inserted and populated"
